$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with uniform run formatting) ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Crime-data grid updates (rows 14-29) ---
# Row 14
$ws.Range("N14").Value = -87.5

# Row 15
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -37.5

# Row 16
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -8.333333333333
$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 37.142857142857
$ws.Range("I16").Value = 78
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 25.806451612903
$ws.Range("L16").Value = 36.842105263157
$ws.Range("M16").Value = 25.806451612903
$ws.Range("N16").Value = -68.032786885245

# Row 17
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 53.333333333333
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = 34.482758620689
$ws.Range("I17").Value = 121
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 14.150943396226
$ws.Range("L17").Value = 53.164556962025
$ws.Range("M17").Value = 152.083333333333
$ws.Range("N17").Value = 23.469387755102

# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 11.764705882352
$ws.Range("L18").Value = -28.301886792452
$ws.Range("M18").Value = 31.034482758620
$ws.Range("N18").Value = -79.459459459459

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 13.725490196078
$ws.Range("I19").Value = 116
$ws.Range("J19").Value = 93
$ws.Range("K19").Value = 24.731182795698
$ws.Range("L19").Value = 41.463414634146
$ws.Range("M19").Value = 157.777777777778
$ws.Range("N19").Value = 45

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -21.739130434782
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -25
$ws.Range("L20").Value = -38.775510204081
$ws.Range("M20").Value = 87.5
$ws.Range("N20").Value = -65.517241379310

# Row 21
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 34.883720930232
$ws.Range("F21").Value = 226
$ws.Range("G21").Value = 189
$ws.Range("H21").Value = 19.576719576719
$ws.Range("I21").Value = 390
$ws.Range("J21").Value = 345
$ws.Range("K21").Value = 13.043478260869
$ws.Range("L21").Value = 20
$ws.Range("M21").Value = 92.118226600985
$ws.Range("N21").Value = -45.682451253481

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -28.571428571428
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = 180

# Row 23
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 125
$ws.Range("F23").Value = 41
$ws.Range("G23").Value = 31
$ws.Range("H23").Value = 32.258064516129
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = 12.727272727272
$ws.Range("L23").Value = 51.219512195122
$ws.Range("M23").Value = 67.567567567567

# Row 24
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -5.714285714285
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = 11.570247933884
$ws.Range("I24").Value = 241
$ws.Range("J24").Value = 218
$ws.Range("K24").Value = 10.550458715596
$ws.Range("L24").Value = 57.516339869281
$ws.Range("M24").Value = 52.531645569620

# Row 25
$ws.Range("C25").Value = 32
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 128.571428571429
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = 69.696969696969
$ws.Range("I25").Value = 174
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 45
$ws.Range("L25").Value = 65.714285714285
$ws.Range("M25").Value = 35.9375

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 100
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -36.363636363636
$ws.Range("L26").Value = 133.333333333333

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = -80
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = -13.333333333333
$ws.Range("L27").Value = 30

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -55.555555555555
$ws.Range("N28").Value = -85.714285714285

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 4
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -42.857142857142
$ws.Range("N29").Value = -81.818181818181

